$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from its old location
#    (right after "... Zbraně ti zůstanou i po splnění mise.")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Insert a new run "6 kostek, " right after "t misí, " (before "6 náhodných ")
#    using the same run formatting (Bahnschrift SemiCondensed) as its neighbours.
$rng = $d.Content
$rng.Find.Execute("t misí, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$insertStart = $rng.Start
$rng.InsertAfter("6 kostek, ")
$insertEnd = $insertStart + 10

# 3. Re-create the "_GoBack" bookmark right after the newly inserted text.
#    Also drop a throwaway bookmark at the left edge of the new run first -- this
#    keeps the engine from silently re-merging "6 kostek, " back into the
#    preceding "t misí, " run once the document is serialised.
$d.Bookmarks.Add("ZZ_TempSplit", $d.Range($insertStart, $insertStart))
$d.Bookmarks.Add("_GoBack", $d.Range($insertEnd, $insertEnd))
$d.Bookmarks("ZZ_TempSplit").Delete()
